$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

$ws.Rows.Item(178).Insert()

$ws.Cells.Item(178, 1).Value = $ws.Cells.Item(179, 1).Value2
$ws.Cells.Item(178, 2).Value = "EoIEECwEC"
$ws.Cells.Item(178, 3).Value = "Elasticity of Industrial Equipment Energy Consumption with respect to Energy Cost"
$ws.Cells.Item(178, 6).Value = "low"
